# Manifold Mass.xlsx update:
#  - add a new "Source" column (G) with manufacturer/spec source notes
#  - update the existing "Notes" column (F) values to reflect specific pipe
#    schedule / fitting callouts instead of the generic manufacturer name
#  - the old generic "weight from another manufacturer..." note moves from
#    column F to the new column G on the two rows that used it
#  - uniform row height across the used range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell: G1 "Source" -----------------------------------------
# Copy the formatting of the existing header cell (F1) onto G1 first so the
# new column matches the header styling (fill/border), then set its text.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G1").Value = "Source"

# --- Notes column (F) rewrites ---------------------------------------------
$ws.Range("F2").Value = "12"" DR17"
$ws.Range("F3").Value = "12""x12""x6"" DR17"
$ws.Range("F4").Value = "12"" DR17"
$ws.Range("F5").ClearContents()
$ws.Range("F6").Value = "12"" "
$ws.Range("F7").Value = "12"" DR17"
$ws.Range("F8").Value = "6"" DR17"
$ws.Range("F9").ClearContents()
$ws.Range("F10").Value = "N60-F PP"

# --- New Source column (G) values ------------------------------------------
# First stamp the formatting of column F (data rows use style index 8) onto
# column G for rows 2-12 so the new cells render the same as their neighbour.
$ws.Range("F2:F12").Copy()
$ws.Range("G2:G12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G2").Value = "EPG, IPS FAB TEE"
$ws.Range("G3").Value = "EPG, IPS FAB BRANCH SADDLE RED. TEE"
$ws.Range("G4").Value = "EPG, IPS FLANGE ADAPTER"
$ws.Range("G5").Value = "weight from another manufacturer since EPG had a 404 error, need their values"
$ws.Range("G6").Value = "Wolesely, IPS HDPE BLIND FLANGE"
$ws.Range("G7").Value = "EPG, PIPE CHART"
$ws.Range("G8").Value = "EPG, IPS FLANGE ADAPTER"
$ws.Range("G9").Value = "weight from another manufacturer since EPG had a 404 error, need their values"
$ws.Range("G10").Value = "Mazzei"
# G11 / G12 stay blank (formatting only, already copied above)

# --- Uniform row heights -----------------------------------------------------
$ws.Range("A1:G12").RowHeight = 15

# --- Column widths (approximate autosize to match new content) -------------
$ws.Columns.Item(1).ColumnWidth = 39
$ws.Columns.Item(6).ColumnWidth = 15.5703125
$ws.Columns.Item(7).ColumnWidth = 72.5703125

# --- Selection, matching the saved cursor position in the source file ------
$ws.Range("G19").Select()
